$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a bare number-looking string need to be forced
# to remain text (matching the source inlineStr cells), otherwise Excel
# auto-converts them to floating point numbers.
$ws.Range('D2').Value = '26.632.06'
$ws.Range('E2').Value = '  +0.96%  '
$ws.Range('D3').Value = '1.631.42'
$ws.Range('E3').Value = '  +0.52%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('E6').Value = '  +2.54%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  +1.51%  '
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('E11').Value = '  +3.43%  '
$ws.Range('D12').Value = '1.859.29'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '1.668.86'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('E14').Value = '  +1.61%  '
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').Value = '26.629.87'
$ws.Range('E16').Value = '  +1.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.14'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.92%  '
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '217.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.51%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.93%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.34'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  +3.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.24%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').Value = '  +0.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('E30').Value = '  -2.84%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  +3.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.98'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.41'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.72%  '
$ws.Range('D36').Value = '1.210.88'
$ws.Range('E36').Value = '  +3.02%  '
$ws.Range('E37').Value = '  +5.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.807'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('D44').Value = '1.770.87'
$ws.Range('E44').Value = '  +0.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.92'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.55'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.70'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.62'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.30%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('E51').Value = '  +0.27%  '
